$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "274.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.71%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.99%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.867"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.02%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06320"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.79%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.887"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.30%"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.58%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.285"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "36.67%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8711"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.89%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1460"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.82%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05052"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.63%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07400"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.53%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-5.41%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09035"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.08%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001571"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.13%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006284"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.16%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005890"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.39%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.453"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.12%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.284"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.10%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.54%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.14%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.899"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.27%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04365"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.09%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001177"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.12%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004264"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.33%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.07%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.29%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04036"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.14%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006672"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.67%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.21%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002089"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.63%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01216"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.29%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005318"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.48%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.381"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.96%"
